# Natmi following Dr Hou advice
#
# The Col4a1-Itgb1 NATMI edge table was regenerated upstream: the ligand-
# and receptor-"expressing cells" counts (cols E and K) went from 1 to 3
# for every row, the average-expression columns (G, M) were recomputed
# against the new data, and every downstream derived column (total
# expression = avg * count, specificity fractions, edge weights) was
# recalculated accordingly. Columns F and L (detection rate) stay 1.
#
# Row order matches the worksheet (row 2 = first data row ... row 17 =
# last data row); each inner array holds columns E..T in order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17)

$data = @(
    @(3, 1, 216.5832213333333, 649.749664, 0.4331411212367192, 0.4331411212367192, 3, 1, 98.946724, 296.840172, 0.2098009692989996, 0.2098009692989996, 21430.20022430025, 192871.8020187022, 0.09087342707871918, 0.09087342707871918),
    @(3, 1, 216.5832213333333, 649.749664, 0.4331411212367192, 0.4331411212367192, 3, 1, 163.0062356666667, 489.018707, 0.345629090707923, 0.3456290907079231, 35304.41561810716, 317739.7405629645, 0.1497061718812575, 0.1497061718812575),
    @(3, 1, 216.5832213333333, 649.749664, 0.4331411212367192, 0.4331411212367192, 3, 1, 65.39610666666668, 196.18832, 0.1386621609326595, 0.1386621609326595, 14163.69944452495, 127473.2950007245, 0.06006028385947852, 0.06006028385947852),
    @(3, 1, 216.5832213333333, 649.749664, 0.4331411212367192, 0.4331411212367192, 3, 1, 144.2727966666667, 432.81839, 0.3059077790604178, 0.3059077790604179, 31247.06705283566, 281223.603475521, 0.132501238417264, 0.132501238417264),
    @(3, 1, 161.954974, 485.864922, 0.3238910133313607, 0.3238910133313606, 3, 1, 98.946724, 296.840172, 0.2098009692989996, 0.2098009692989996, 16024.91411280518, 144224.2270152466, 0.06795264854415466, 0.06795264854415466),
    @(3, 1, 161.954974, 485.864922, 0.3238910133313607, 0.3238910133313606, 3, 1, 163.0062356666667, 489.018707, 0.345629090707923, 0.3456290907079231, 26399.67065923287, 237597.0359330958, 0.111946156426186, 0.111946156426186),
    @(3, 1, 161.954974, 485.864922, 0.3238910133313607, 0.3238910133313606, 3, 1, 65.39610666666668, 196.18832, 0.1386621609326595, 0.1386621609326595, 10591.22475490123, 95321.02279411105, 0.04491142781519528, 0.04491142781519528),
    @(3, 1, 161.954974, 485.864922, 0.3238910133313607, 0.3238910133313606, 3, 1, 144.2727966666667, 432.81839, 0.3059077790604178, 0.3059077790604179, 23365.69703305728, 210291.2732975156, 0.09908078054582473, 0.09908078054582473),
    @(3, 1, 0.4608033333333333, 1.38241, 0.0009215527926904059, 0.0009215527926904059, 3, 1, 98.946724, 296.840172, 0.2098009692989996, 0.2098009692989996, 45.59498024161333, 410.35482217452, 0.0001933426691666472, 0.0001933426691666472),
    @(3, 1, 0.4608033333333333, 1.38241, 0.0009215527926904059, 0.0009215527926904059, 3, 1, 163.0062356666667, 489.018707, 0.345629090707923, 0.3456290907079231, 75.11381674931889, 676.02435074387, 0.0003185154537769321, 0.0003185154537769321),
    @(3, 1, 0.4608033333333333, 1.38241, 0.0009215527926904059, 0.0009215527926904059, 3, 1, 65.39610666666668, 196.18832, 0.1386621609326595, 0.1386621609326595, 30.13474393902223, 271.2126954512, 0.0001277845016479788, 0.0001277845016479788),
    @(3, 1, 0.4608033333333333, 1.38241, 0.0009215527926904059, 0.0009215527926904059, 3, 1, 144.2727966666667, 432.81839, 0.3059077790604178, 0.3059077790604179, 66.48138561332222, 598.3324705199, 0.0002819101680988477, 0.0002819101680988478),
    @(3, 1, 121.0302313333333, 363.090694, 0.2420463126392298, 0.2420463126392298, 3, 1, 98.946724, 296.840172, 0.2098009692989996, 0.2098009692989996, 11975.54489539549, 107779.9040585594, 0.0507815510069591, 0.0507815510069591),
    @(3, 1, 121.0302313333333, 363.090694, 0.2420463126392298, 0.2420463126392298, 3, 1, 163.0062356666667, 489.018707, 0.345629090707923, 0.3456290907079231, 19728.68241151252, 177558.1417036127, 0.08365824694670265, 0.08365824694670265),
    @(3, 1, 121.0302313333333, 363.090694, 0.2420463126392298, 0.2420463126392298, 3, 1, 65.39610666666668, 196.18832, 0.1386621609326595, 0.1386621609326595, 7914.90591816601, 71234.1532634941, 0.03356266475633769, 0.03356266475633769),
    @(3, 1, 121.0302313333333, 363.090694, 0.2420463126392298, 0.2420463126392298, 3, 1, 144.2727966666667, 432.81839, 0.3059077790604178, 0.3059077790604179, 17461.36995567363, 157152.3296010627, 0.07404384992923033, 0.07404384992923033)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $vals = $data[$i]
    $col = 5
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}
